$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$ws = $wb.Worksheets.Item("About")

# Remove the date stamp in C1 entirely (value + formatting)
$ws.Range("C1").Clear()

# Remove the hyperlink that lived on B6 (it will be removed along with rows 4:7 below)
$ws.Hyperlinks.Delete()

# Remove the old source detail rows (year, paper title, url, page) - rows 4 through 7
$ws.Rows("4:7").Delete()

# Replace the (now single) notes line with the new explanatory text
$ws.Range("A9").Value = "In the US, we set this to 0 so that increasing EV chargers does not induce additional deployment."

# Update the "Source:" value cell
$ws.Range("B3").Value = "None"

# --- "EoCSoEVMS" sheet ---
$ws2 = $wb.Worksheets.Item("EoCSoEVMS")
$ws2.Range("B2").Value = 0

# Match the selection state left behind in each sheet
$ws2.Activate() | Out-Null
$ws2.Range("B3").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("A4:XFD7").Select() | Out-Null
